# LIMS Update from Lakshman System
#
# 1. Insert a new blank worksheet "Sheet3" right after "StabilityRole"
#    (and before "MCRole").
# 2. Rename "Plant-3" (and its "Plant-3-QC/-PF/-QA" variants) to "Plant-2"
#    on the UserCreation sheet.
# 3. Update TestDetails' remembered selection (B16), no longer the active tab.
# 4. Append a new worksheet "RMTestDetails" at the end with the RM
#    Qualitative/Quantitative test rows, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Insert empty "Sheet3" after StabilityRole, before MCRole ---------
$stabilityRole = $wb.Worksheets.Item("StabilityRole")
$sheet3 = $wb.Worksheets.Add($null, $stabilityRole)
$sheet3.Name = "Sheet3"

# --- 2. UserCreation: Plant-3 -> Plant-2 (and QC/PF/QA variants) ---------
$userCreation = $wb.Worksheets.Item("UserCreation")
for ($r = 2; $r -le 9; $r++) {
    $locCell = $userCreation.Cells.Item($r, 2)
    $locVal = [string]$locCell.Value()
    $locCell.Value = $locVal.Replace("Plant-3", "Plant-2")

    $deptCell = $userCreation.Cells.Item($r, 6)
    $deptVal = [string]$deptCell.Value()
    $deptCell.Value = $deptVal.Replace("Plant-3", "Plant-2")
}
$userCreation.Range("D14").Select()

# --- 3. TestDetails: selection moves to B16 (no longer the active tab) ---
$testDetails = $wb.Worksheets.Item("TestDetails")
$testDetails.Range("B16").Select()

# --- 4. Append "RMTestDetails" sheet at the end ---------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$rmTestDetails = $wb.Worksheets.Add($null, $lastSheet)
$rmTestDetails.Name = "RMTestDetails"

$rmTestDetails.Cells.Item(1, 1).Value = "TestName"
$rmTestDetails.Cells.Item(1, 2).Value = "TestType"

$rmData = @(
    @("RM Qualitative Test -2", "Qualitative"),
    @("RM Qualitative Test -1", "Qualitative"),
    @("RM Qualitative Test -3", "Qualitative"),
    @("RM Qualitative Test -4", "Qualitative"),
    @("RM Qualitative Test -5", "Qualitative"),
    @("RM Qualitative Test -6", "Qualitative"),
    @("RM Quantitative Test-1", "Quantitative"),
    @("RM Quantitative Test-2", "Quantitative"),
    @("RM Quantitative Test-3", "Quantitative"),
    @("RM Quantitative Test-4", "Quantitative"),
    @("RM Quantitative Test-5", "Quantitative"),
    @("RM Quantitative Test-6", "Quantitative")
)

# Row 3 (RM Qualitative Test -2) is entered before row 2 (RM Qualitative
# Test -1) so that the shared-string table order matches the source
# workbook (string "-2" registered ahead of "-1").
$rmTestDetails.Cells.Item(3, 1).Value = $rmData[0][0]
$rmTestDetails.Cells.Item(3, 2).Value = $rmData[0][1]
$rmTestDetails.Cells.Item(2, 1).Value = $rmData[1][0]
$rmTestDetails.Cells.Item(2, 2).Value = $rmData[1][1]
for ($i = 2; $i -lt $rmData.Length; $i++) {
    $row = $i + 2
    $rmTestDetails.Cells.Item($row, 1).Value = $rmData[$i][0]
    $rmTestDetails.Cells.Item($row, 2).Value = $rmData[$i][1]
}

$rmTestDetails.Range("A16").Select()
$rmTestDetails.Activate()
